$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("Accountability Documentation", "Data consumer", "References or bibliography to related work for more context"),
  @("Accountability Documentation", "Both-Data generator", "data transparency essentials, demographics, instrumentation and methods for collection and processing."),
  @("Accountability Documentation", "Both-Data consumer", "Depending on the nature of the data, one should try to ensure accuracy, avoid missing entries, biases, ensure details and patients' data and response information accurately"),
  @("Accountability Documentation", "Both-Data generator", "- Known Biases and Bias mitigation strategies`n- Ethical compliance testing (probably in Testing Apendix C)"),
  @("Accountability Documentation", "Both-Data consumer", "Social Determinants of Health."),
  @("Data Card", "Both-Data generator", "It is notable that this is case by case.  Work with protected populations as children or the Indigenous have special protections.  In general, Social Determinants of Health, and environmental data are additional important information."),
  @("Data Card", "Both-Data consumer", "- Generalizability and Limitations: Discuss the generalizability of the data to different populations/settings and any limitations that should be considered by those developing models with the data."),
  @("Data Card", "Both-Data generator", "Some potential extra info could be:`n- Informed Consent: (or IRB approvals etc.)`n- Potential bias sources`n- Ethical Review: whether the dataset has undergone any ethical review?"),
  @("Data Card", "Data generator", "n/s"),
  @("Data Card", "Both-Data generator", "preferable way to cite it (i.e published paper)"),
  @("Data Card", "Data generator", "degree of ""missingness"" (and reasons for) in the data. "),
  @("Dataset Nutrition Label", "Both-Data generator", "Used metrics and what are thier references"),
  @("Dataset Nutrition Label", "Data consumer", "A list of reasons for why data may not have been collected or is missing (i.e. exclusion criteria, invalid results, etc) could be very informative when performing analyses of a dataset. "),
  @("Dataset Nutrition Label", "Both-Data consumer", "1. Details on any methods used to identify and mitigate biases in the dataset, which is critical for ensuring fair AI/ML model development.`n2. Specific breakdowns of demographic variables such as age, gender, race, and socioeconomic status to understand the population represented in the dataset.`n3. Documentation of ethical approvals, informed consent, and measures taken to protect patient privacy.`n"),
  @("Dataset Nutrition Label", "Both-Data generator", "1. Information on cross-validation methods and results to assess model robustness. [Assessment of how findings can be generalized to broader populations.] Specifics on clinical variables collected, including measurement units and methods.`n2. Inclusion of environmental, social, or economic factors that may influence the data.`n Consideration of historical or temporal events that could impact the dataset (e.g., pandemics, policy changes)."),
  @("Dataset Nutrition Label", "Both-Data consumer", "Data distrubution for each metric "),
  @("Datasheet", "Both-Data generator", "IRB protocols consent from subjects is critical and must be entered. Data entries and database maintenance is also very important. One should add data ethnicity, gender and age details besides the main data readings  such as EKG, retina scan etc. "),
  @("Datasheet", "Data consumer", "Dataset Versioning, Data Bias and Fairness"),
  @("Datasheet", "Both-Data consumer", "human gender, age, ethnicity and other background information is very useful to assess their data for diagnosis. For example, do they have parental history of diabetes if they are also labeled diabetic? All this would be important information to have in addition to their EEG, EKG and retina scan readings etc. "),
  @("Datasheet", "Both-Data consumer", "- Model-specific considerations: any guidelines or warnings about particular AI/ML models that may amplify biases or otherwise misinterpret the data. etc."),
  @("Datasheet", "Both-Data generator", "- More broader demographic details beside subpopulations to allow for a deeper understanding of representation and potential biases.`n- Dataset usage history and feedback`n- User reporting mechanism: a method for users to report any issues they discover within the dataset, e.g., biases, errors, etc.`n- Model-specific considerations: any guidelines or warnings about particular AI/ML models that may amplify biases or otherwise misinterpret the data. etc."),
  @("Datasheet", "Data consumer", "Sampling process"),
  @("Healthsheet", "Both-Data generator", "Missingness of data"),
  @("Healthsheet", "Both-Data consumer", "Misssingness of data. Many factors can influence whether data is missing from an individual case (e.g., participant). It would be good to know why it was missing."),
  @("Healthsheet", "Both-Data generator", "1. Information on how potential biases in the dataset were identified and mitigated during data collection and processing.`n2.  Details on ethical approvals, informed consent, and measures to protect patient privacy.`n3.  Clear documentation of the data's origin, including any transformations or processing steps.`n4. Guidance on the appropriate and inappropriate use of the dataset for developing AI/ML models, including known limitations."),
  @("Healthsheet", "Both-Data consumer", "1. Explicit details on any biases identified within the dataset and the steps taken to mitigate them.`n2. Detailed tracking of data lineage, including sources, transformations, and preprocessing steps.`n"),
  @("Healthsheet", "Data generator", "-Data security: for example, who might be affected by any possible data leaks, what segments of data must have a greater security factor added to it, possible data breach factors, and how to prevent them? `n-Priority ranking: How meaningful is each dataset to the data as a whole?")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
